# Aula T2-S01 e T2-S02, updates de código
#
# Moves the main content placeholder ("Rectangle 23" / "Rectangle 3", the
# 2nd shape on each slide) up slightly: xfrm offset y goes from its
# previous value (913284 / 900000 EMU) to 841276 EMU on slides 2-21.
# x (Left) and the shape's cx/cy (Width/Height) are left untouched.
#
# PowerPoint's COM object model works in points (1 pt = 12700 EMU), and
# internally rounds the point value through a single-precision float
# before converting back to EMU, so the literal below (66.24221 pt) is
# the value that round-trips to exactly 841276 EMU.

$p = $ppt.ActivePresentation

$targetTopPt = 66.24221

for ($i = 2; $i -le 21; $i++) {
    $slide = $p.Slides.Item($i)
    $shape = $slide.Shapes.Item(2)
    $shape.Top = $targetTopPt
}
